$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 / Row 13: swap Polkadot <-> WrappedEther (B, C), then update D, E ---
$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value = '1.647.28'
$ws.Cells.Item(12, 5).Value = '  -0.53%  '

$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.253'
$ws.Cells.Item(13, 5).Value = '  -1.21%  '

# --- Remaining rows: update Price (D) and Volume(1h) (E) columns ---
$ws.Cells.Item(2, 4).Value = '25.880.81'
$ws.Cells.Item(2, 5).Value = '  -1.31%  '
$ws.Cells.Item(3, 4).Value = '1.641.44'
$ws.Cells.Item(3, 5).Value = '  -0.85%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.006'
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '215.79'
$ws.Cells.Item(5, 5).Value = '  -0.28%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.5030'
$ws.Cells.Item(6, 5).Value = '  -2.02%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.006'
$ws.Cells.Item(7, 5).Value = '  -0.17%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2570'
$ws.Cells.Item(8, 5).Value = '  -1.20%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06373'
$ws.Cells.Item(9, 5).Value = '  -1.37%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.62'
$ws.Cells.Item(10, 5).Value = '  -1.58%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07756'
$ws.Cells.Item(11, 5).Value = '  -1.06%  '
$ws.Cells.Item(14, 4).Value = '1.865.19'
$ws.Cells.Item(14, 5).Value = '  -1.06%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.5455'
$ws.Cells.Item(15, 5).Value = '  -1.55%  '
$ws.Cells.Item(16, 4).Value = '0.0₅7882'
$ws.Cells.Item(16, 5).Value = '  -2.05%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '64.05'
$ws.Cells.Item(17, 5).Value = '  -0.27%  '
$ws.Cells.Item(18, 4).Value = '25.927.34'
$ws.Cells.Item(18, 5).Value = '  -1.19%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.005'
$ws.Cells.Item(19, 5).Value = '  -0.16%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '202.06'
$ws.Cells.Item(20, 5).Value = '  -4.20%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.386'
$ws.Cells.Item(21, 5).Value = '  -0.74%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '9.894'
$ws.Cells.Item(22, 5).Value = '  -1.98%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.973'
$ws.Cells.Item(23, 5).Value = '  -0.95%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '1.006'
$ws.Cells.Item(24, 5).Value = '  -0.11%  '
$ws.Cells.Item(25, 5).Value = '  +6.46%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '140.73'
$ws.Cells.Item(26, 5).Value = '  -2.83%  '
$ws.Cells.Item(27, 5).Value = '  -3.89%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '15.61'
$ws.Cells.Item(28, 5).Value = '  -1.64%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '6.770'
$ws.Cells.Item(29, 5).Value = '  -3.59%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.243'
$ws.Cells.Item(30, 5).Value = '  -0.02%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.04960'
$ws.Cells.Item(31, 5).Value = '  -2.83%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.259'
$ws.Cells.Item(32, 5).Value = '  -3.17%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.192'
$ws.Cells.Item(33, 5).Value = '  -1.26%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.546'
$ws.Cells.Item(34, 5).Value = '  -1.26%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.374'
$ws.Cells.Item(35, 5).Value = '  +0.84%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.634'
$ws.Cells.Item(36, 5).Value = '  -3.70%  '
$ws.Cells.Item(37, 5).Value = '  -3.66%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.5627'
$ws.Cells.Item(38, 5).Value = '  -1.97%  '
$ws.Cells.Item(39, 4).Value = '1.148.29'
$ws.Cells.Item(39, 5).Value = '  -1.59%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.01566'
$ws.Cells.Item(40, 5).Value = '  -1.63%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.006'
$ws.Cells.Item(41, 5).Value = '  -0.08%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.674'
$ws.Cells.Item(42, 5).Value = '  -0.78%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '99.94'
$ws.Cells.Item(43, 5).Value = '  -0.44%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.8065'
$ws.Cells.Item(44, 5).Value = '  -2.26%  '
$ws.Cells.Item(45, 4).Value = '1.776.31'
$ws.Cells.Item(45, 5).Value = '  -1.14%  '
$ws.Cells.Item(46, 5).Value = '  +0.14%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.4548'
$ws.Cells.Item(47, 5).Value = '  +0.02%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.007'
$ws.Cells.Item(48, 5).Value = '  -0.01%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '54.83'
$ws.Cells.Item(49, 5).Value = '  -1.15%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.05058'
$ws.Cells.Item(50, 5).Value = '  -0.41%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.004'
$ws.Cells.Item(51, 5).Value = '  -0.31%  '
